# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country names (rows whose A-column label changes because the
# underlying shared-string position moved) and refresh their stats ---

# Canada <-> Portugal swap (rows 18/19)
$ws.Range("A18").Value = "Canada"
$ws.Range("B18").Value = 4610
$ws.Range("C18").Value = 567
$ws.Range("D18").Value = 228
$ws.Range("E18").Value = 4343
$ws.Range("F18").Value = 120
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 39

$ws.Range("A19").Value = "Portugal"
$ws.Range("B19").Value = 4268
$ws.Range("C19").Value = 724
$ws.Range("D19").Value = 43
$ws.Range("E19").Value = 4149
$ws.Range("F19").Value = 71
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 76

# Chequia <-> Malasia swap (rows 25/26)
$ws.Range("A25").Value = "Chequia"
$ws.Range("B25").Value = 2279
$ws.Range("C25").Value = 354
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 2259
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 9

$ws.Range("A26").Value = "Malasia"
$ws.Range("B26").Value = 2161
$ws.Range("C26").Value = 130
$ws.Range("D26").Value = 259
$ws.Range("E26").Value = 1876
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 26

# Monaco / Puerto Rico / Macao reorder (rows 126/127/128)
$ws.Range("A126").Value = "Monaco"
$ws.Range("B126").Value = 42
$ws.Range("C126").Value = 9
$ws.Range("D126").Value = 1
$ws.Range("E126").Value = 41
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 0

$ws.Range("A127").Value = "Puerto Rico"
$ws.Range("B127").Value = 39
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 1
$ws.Range("E127").Value = 36
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 2

$ws.Range("A128").Value = "Macao"
$ws.Range("B128").Value = 34
$ws.Range("C128").Value = 1
$ws.Range("D128").Value = 10
$ws.Range("E128").Value = 24
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 0

# --- Plain statistic refreshes (no row reordering) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 94425
$ws.Range("C4").Value = 8990
$ws.Range("D4").Value = 2447
$ws.Range("E4").Value = 90549
$ws.Range("F4").Value = 2463
$ws.Range("G4").Value = 134
$ws.Range("H4").Value = 1429

# Suiza (row 12)
$ws.Range("B12").Value = 12928
$ws.Range("C12").Value = 1117
$ws.Range("E12").Value = 11800
$ws.Range("G12").Value = 39
$ws.Range("H12").Value = 231

# Austria (row 15)
$ws.Range("B15").Value = 7610
$ws.Range("C15").Value = 701
$ws.Range("E15").Value = 7327

# Pakistan (row 34)
$ws.Range("B34").Value = 1331
$ws.Range("C34").Value = 130
$ws.Range("E34").Value = 1299

# India (row 44)
$ws.Range("B44").Value = 887
$ws.Range("C44").Value = 160
$ws.Range("E44").Value = 794

# Serbia (row 57)
$ws.Range("D57").Value = 42
$ws.Range("E57").Value = 478

# Lituania (row 66)
$ws.Range("B66").Value = 358
$ws.Range("C66").Value = 59
$ws.Range("E66").Value = 352

# Marruecos (row 67)
$ws.Range("D67").Value = 11
$ws.Range("E67").Value = 301
$ws.Range("G67").Value = 10
$ws.Range("H67").Value = 21

# --- Update the "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 18:44"
